$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The user fills in the last empty "Value" placeholder (previously shown as "-")
# in the converter table with an actual decimal value (72). This causes the
# dependent formulas in F12 (hex), G12 (75% opacity) and H12 (50% opacity) to
# recalculate from errors into real results.
$ws.Range("E12").Value = 72

# After typing the value and pressing Enter, the active selection moves down
# to the next cell.
$ws.Range("E13").Select()
